# "Added 6 Test cases" — set the Sheet1 "Login" row's password column (B1)
# to "Pass" (this also introduces a new shared string "Pass").
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

$ws1.Activate()
$ws1.Range("B1").Value = "Pass"

# Reflect the resulting selection (B1:C1, same as the authored workbook).
$ws1.Range("B1:C1").Select()
